$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.959.57"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "2.929.92"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.02%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.57"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.45%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.79"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.42%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("E8").Value = "  +0.01%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.634"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.07%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.51"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -1.57%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0881"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +2.06%  "

$ws.Range("E12").Value = "  +0.91%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.74"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -2.00%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.91"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").Value = "3.393.47"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "2.932.97"
$ws.Range("E16").Value = "  +0.31%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.986"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "51.970.91"
$ws.Range("E18").Value = "  -0.89%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.37"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.84%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -1.76%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -4.77%  "

$ws.Range("E22").Value = "  +0.25%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.07"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.14%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.53"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("E25").Value = "  +0.36%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.188"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +14.28%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.11"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.49%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +15.81%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +13.66%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.65"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.19%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.18"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("E33").Value = "  +1.88%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.06"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -1.80%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.30"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -1.52%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0445"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("E37").Value = "  +0.06%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -2.05%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.50"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -1.79%  "

$ws.Range("E40").Value = "  -3.03%  "

$ws.Range("E41").Value = "  +1.31%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +2.82%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.98"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -5.72%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.67"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -2.67%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.52%  "

$ws.Range("E46").Value = "  -2.08%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.48"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").Value = "2.137.41"
$ws.Range("E48").Value = "  -3.74%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.248"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -7.28%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0333"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.41%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.19"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -0.14%  "
